# Refactor the synthetic "statut" color array:
#   black square (noir)  -> blue book  (bleu)
#   orange square (orange) stays "orange" label but icon -> orange book
#   red square (rouge) stays "rouge" label but icon -> red book
#
# Column A ("statut") holds the emoji square, column B ("statut_label")
# holds the French color name. Only the "black/noir" entries change
# their label text (-> "bleu"); orange/rouge keep their existing labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the used range so we don't hard-code row numbers.
$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $statut = $ws.Cells.Item($r, 1).Value2

    if ($statut -eq "⬛") {
        $ws.Cells.Item($r, 1).Value = "📘"
        $ws.Cells.Item($r, 2).Value = "bleu"
    }
    elseif ($statut -eq "🟧") {
        $ws.Cells.Item($r, 1).Value = "📙"
    }
    elseif ($statut -eq "🟥") {
        $ws.Cells.Item($r, 1).Value = "📕"
    }
}
